$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 6 data (Day 5 entry)
$ws.Range("A6").Value = "Day 5"

# Copy the date cell's format (so it reuses the existing date style) then set the value
$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B6").Value2 = 45807

$ws.Range("C6").Value = "Product of Array Except Self"
$ws.Range("D6").Value = "Container With Most Water"
$ws.Range("E6").Value = "Min Stack"
$ws.Range("F6").Value = "Prefix-Suffix, Two Pointer, Stack Design"
$ws.Range("G6").Value = "S"
$ws.Range("H6").Value = "YES"

# Update selection to H6, matching the saved view state
$ws.Range("H6").Select()
